$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "Globo"
$ws.Range("B34").Value = "RJ TV 2"
$ws.Range("C34").Value = "Esportes"
$ws.Range("D34").Value = "2025-04-01T19:38"
$ws.Range("E34").Value = "Neutro"
$ws.Range("F34").Value = "Mudança na presidência do Americano. Tolentino Reis é destituído por votação do Conselho Deliberativo. Laila Póvoa assume. *nota coberta*"
